# A new weekly price record (row 39, 2022-03-17) is inserted before the
# existing row 39 (2021-03-09), pushing all rows from 39..47 down to 40..48.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(39).Insert()

$ws.Cells.Item(39,1).Value = 4
$ws.Cells.Item(39,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(39,3).Value = "Los Lagos"
$ws.Cells.Item(39,4).Value = 44637
$ws.Cells.Item(39,5).Value = 10
$ws.Cells.Item(39,6).Value = 100112043
$ws.Cells.Item(39,7).Value = "Pepino dulce"
$ws.Cells.Item(39,8).Value = "Cultivar IV Región"
$ws.Cells.Item(39,9).Value = "Especial"
$ws.Cells.Item(39,10).Value = 50
$ws.Cells.Item(39,11).Value = 21000
$ws.Cells.Item(39,12).Value = 21000
$ws.Cells.Item(39,13).Value = 21000
$ws.Cells.Item(39,14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(39,15).Value = "Provincia de Limarí"
$ws.Cells.Item(39,16).Value = 1167
$ws.Cells.Item(39,17).Value = 18
$ws.Cells.Item(39,18).Value = "Hortaliza"
